$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46/47: Coin name and Link swap (Cronos <-> Algorand) plus refreshed price data
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"

# Refresh Price (D) and Volume(1h) (E) columns for all rows.
# D values that parse as plain numbers need a leading apostrophe so Excel
# keeps them as literal text (matching the source feed's inlineStr cells)
# instead of silently converting them to numeric cells.
$ws.Range("D2").Value = "25.028.77"
$ws.Range("E2").Value = "  -3.60%  "
$ws.Range("D3").Value = "1.640.55"
$ws.Range("E3").Value = "  -5.60%  "
$ws.Range("D4").Value = "'0.9983"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'232.42"
$ws.Range("E5").Value = "  -5.91%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.4719"
$ws.Range("E7").Value = "  -6.34%  "
$ws.Range("D8").Value = "'0.2543"
$ws.Range("E8").Value = "  -6.77%  "
$ws.Range("D9").Value = "'0.06062"
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("D10").Value = "'0.06999"
$ws.Range("E10").Value = "  -3.64%  "
$ws.Range("D11").Value = "1.637.62"
$ws.Range("E11").Value = "  -5.76%  "
$ws.Range("D12").Value = "'14.28"
$ws.Range("E12").Value = "  -6.41%  "
$ws.Range("D13").Value = "'4.293"
$ws.Range("E13").Value = "  -9.71%  "
$ws.Range("D14").Value = "'0.5649"
$ws.Range("E14").Value = "  -14.04%  "
$ws.Range("D15").Value = "'73.14"
$ws.Range("E15").Value = "  -5.81%  "
$ws.Range("D16").Value = "'0.9998"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'0.9991"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "25.004.60"
$ws.Range("E18").Value = "  -3.75%  "
$ws.Range("D19").Value = "'11.23"
$ws.Range("E19").Value = "  -5.52%  "
$ws.Range("D20").Value = "'0.000006558"
$ws.Range("E20").Value = "  -4.02%  "
$ws.Range("D21").Value = "1.846.35"
$ws.Range("E21").Value = "  -5.85%  "
$ws.Range("D22").Value = "'4.260"
$ws.Range("E22").Value = "  -7.81%  "
$ws.Range("D23").Value = "'8.439"
$ws.Range("E23").Value = "  -4.23%  "
$ws.Range("D24").Value = "'5.189"
$ws.Range("E24").Value = "  -4.80%  "
$ws.Range("D25").Value = "'132.47"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("D26").Value = "'14.82"
$ws.Range("E26").Value = "  -2.95%  "
$ws.Range("D27").Value = "'1.367"
$ws.Range("E27").Value = "  -6.07%  "
$ws.Range("D28").Value = "'103.41"
$ws.Range("E28").Value = "  -1.93%  "
$ws.Range("D29").Value = "'1.629"
$ws.Range("E29").Value = "  -9.10%  "
$ws.Range("D30").Value = "'3.866"
$ws.Range("E30").Value = "  -3.18%  "
$ws.Range("D31").Value = "'0.07540"
$ws.Range("E31").Value = "  -7.35%  "
$ws.Range("D32").Value = "'3.489"
$ws.Range("E32").Value = "  -6.27%  "
$ws.Range("D33").Value = "'0.9994"
$ws.Range("D34").Value = "'0.04213"
$ws.Range("E34").Value = "  -11.62%  "
$ws.Range("D35").Value = "'2.568"
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("D36").Value = "'0.9289"
$ws.Range("E36").Value = "  -6.91%  "
$ws.Range("D37").Value = "'0.5861"
$ws.Range("E37").Value = "  -4.27%  "
$ws.Range("D38").Value = "'2.567"
$ws.Range("E38").Value = "  -6.24%  "
$ws.Range("D39").Value = "'0.8619"
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("D40").Value = "'0.9996"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'0.01475"
$ws.Range("E41").Value = "  -8.60%  "
$ws.Range("D42").Value = "'97.69"
$ws.Range("E42").Value = "  -3.12%  "
$ws.Range("D43").Value = "'1.755"
$ws.Range("E43").Value = "  -10.29%  "
$ws.Range("D44").Value = "'0.3644"
$ws.Range("E44").Value = "  -7.56%  "
$ws.Range("D45").Value = "'4.619"
$ws.Range("E45").Value = "  -8.14%  "
$ws.Range("D46").Value = "'0.1087"
$ws.Range("E46").Value = "  -8.13%  "
$ws.Range("D47").Value = "'0.05190"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("D48").Value = "'6.042"
$ws.Range("E48").Value = "  -5.14%  "
$ws.Range("D50").Value = "'0.9990"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "'28.55"
$ws.Range("E51").Value = "  -7.77%  "
